# feat: add 2022-Q3 data
#
# - Insert a new "2022-Q3" sheet (a copy of the "2022-Q2" sheet's layout)
#   positioned right after "总计" and before "2022-Q2", populated with the
#   new quarter's figures.
# - Update the "总计" (totals) summary sheet: push the existing quarter
#   rows down by one and add a new top row for 2022-Q3.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the "2022-Q2" sheet to create "2022-Q3", placed before it.
# ------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Fill in the new quarter's numbers (row 2 = fund A, row 3 = fund C).
# These cells hold text-like numeric strings in the source data, so force
# the number format to Text first (mirrors typing '0.26 into Excel).
$q3Sheet.Range("D2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "0.26"
$q3Sheet.Range("E2").NumberFormat = "@"
$q3Sheet.Range("E2").Value = "86.78"
$q3Sheet.Range("F2").NumberFormat = "@"
$q3Sheet.Range("F2").Value = "5.07"
$q3Sheet.Range("G2").NumberFormat = "@"
$q3Sheet.Range("G2").Value = "0.0132"
$q3Sheet.Range("H2").Value = 10

$q3Sheet.Range("D3").NumberFormat = "@"
$q3Sheet.Range("D3").Value = "0.20"
$q3Sheet.Range("E3").NumberFormat = "@"
$q3Sheet.Range("E3").Value = "86.78"
$q3Sheet.Range("F3").NumberFormat = "@"
$q3Sheet.Range("F3").Value = "5.07"
$q3Sheet.Range("G3").NumberFormat = "@"
$q3Sheet.Range("G3").Value = "0.0101"
$q3Sheet.Range("H3").Value = 10

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: shift existing rows down by one and
#    insert the 2022-Q3 row at the top of the data.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Push rows 2-5 down to rows 3-6, working from the bottom up so nothing
# is overwritten before it's been read.
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(6, 3).Value = 4
$totalSheet.Cells.Item(6, 4).Value = 0.57
$totalSheet.Cells.Item(5, 1).Copy()
$totalSheet.Cells.Item(6, 1).PasteSpecial(-4122)

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(5, 3).Value = 3
$totalSheet.Cells.Item(5, 4).Value = 0.16

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(4, 3).Value = 2
$totalSheet.Cells.Item(4, 4).Value = 0.04

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q2"
$totalSheet.Cells.Item(3, 3).Value = 2
$totalSheet.Cells.Item(3, 4).Value = 0.05

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.02
